# Commit: "Adjusted Excel read file to properly test line skips"
#
# SecondSheet's data block originally started at row 3 (header) / row 4
# (first data row) and ran through row 103. To give the krangl reader test
# some leading blank rows to skip over, insert two blank rows above the
# existing data block, pushing the header down to row 5 and the last data
# row down to row 105.
#
# This also updates the active sheet (SecondSheet becomes the active tab
# instead of ThirdSheet) and the selected cell on SecondSheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SecondSheet")

# Insert two new blank rows above the current row 3, shifting all existing
# data (rows 3:103) down to rows 5:105.
$ws.Rows("3:4").Insert()

# Make SecondSheet the active sheet/tab (was ThirdSheet before the edit).
$ws.Activate()

# Update the visible selection on SecondSheet to E6 (previously F8).
$ws.Range("E6").Select()
